$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5271
$ws.Range("J18").Value = 300
$ws.Range("L18").Value = 300
$ws.Range("N18").Value = -868
$ws.Range("H70").Value = 63727764
$ws.Range("I70").Value = 50002580
$ws.Range("J70").Value = 69446590
$ws.Range("K70").Value = 150007740
$ws.Range("L70").Value = 208339770
$ws.Range("M70").Value = -150007470
$ws.Range("N70").Value = -208340310
$ws.Range("H73").Value = 63727764
$ws.Range("I73").Value = 50002580
$ws.Range("J73").Value = 69446590
$ws.Range("K73").Value = 150007740
$ws.Range("L73").Value = 208339770
$ws.Range("M73").Value = -150006804
$ws.Range("N73").Value = -208341642
$ws.Range("H74").Value = 50009148
$ws.Range("I74").Value = 71435000
$ws.Range("J74").Value = 15501.333
$ws.Range("K74").Value = 71435000
$ws.Range("L74").Value = 15501.333
$ws.Range("M74").Value = -71434064
$ws.Range("N74").Value = -17373.333
$ws.Range("H77").Value = 50009148
$ws.Range("I77").Value = 71435000
$ws.Range("J77").Value = 15501.333
$ws.Range("K77").Value = 357175000
$ws.Range("L77").Value = 77506.66500000001
$ws.Range("M77").Value = -357170320
$ws.Range("N77").Value = -86866.66500000001
$ws.Range("H100").Value = 2362.5
$ws.Range("I100").Value = 1781.091
$ws.Range("K100").Value = 1781.091
$ws.Range("M100").Value = -1240.091
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1571.6
$ws.Range("I102").Value = 1674
$ws.Range("J102").Value = 1332.6666
$ws.Range("K102").Value = 1674
$ws.Range("L102").Value = 1332.6666
$ws.Range("M102").Value = -52
$ws.Range("N102").Value = -4576.6666
$ws.Range("H122").Value = 14305.5
$ws.Range("I122").Value = 16860.385
$ws.Range("K122").Value = 50581.155
$ws.Range("M122").Value = -48131.155
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 94956
$ws.Range("J132").Value = 94956
$ws.Range("L132").Value = 94956
$ws.Range("N132").Value = -105076
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10735.167
$ws.Range("I31").Value = 4011.7
$ws.Range("J31").Value = 14096.9
$ws.Range("K31").Value = 4011.7
$ws.Range("L31").Value = 14096.9
$ws.Range("M31").Value = -3716.7
$ws.Range("N31").Value = -14686.9
$ws.Range("H34").Value = 10735.167
$ws.Range("I34").Value = 4011.7
$ws.Range("J34").Value = 14096.9
$ws.Range("K34").Value = 4011.7
$ws.Range("L34").Value = 14096.9
$ws.Range("M34").Value = -3809.7
$ws.Range("N34").Value = -14500.9
$ws.Range("H51").Value = 100000
$ws.Range("J51").Value = 100000
$ws.Range("L51").Value = 100000
$ws.Range("N51").Value = -101472
$ws.Range("H58").Value = 20842376
$ws.Range("I58").Value = 71430820
$ws.Range("J58").Value = 11843.471
$ws.Range("K58").Value = 71430820
$ws.Range("L58").Value = 11843.471
$ws.Range("M58").Value = -71430617
$ws.Range("N58").Value = -12249.471
$ws.Range("H59").Value = 72222.22
$ws.Range("J59").Value = 93333.336
$ws.Range("L59").Value = 93333.336
$ws.Range("N59").Value = -95623.336
$ws.Range("H60").Value = 22797.8
$ws.Range("I60").Value = 11333.333
$ws.Range("J60").Value = 39994.5
$ws.Range("K60").Value = 11333.333
$ws.Range("L60").Value = 39994.5
$ws.Range("M60").Value = -10822.333
$ws.Range("N60").Value = -41016.5
$ws.Range("H61").Value = 100000
$ws.Range("J61").Value = 100000
$ws.Range("L61").Value = 100000
$ws.Range("N61").Value = -100696
$ws.Range("H74").Value = 1000000000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1000000000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 8424.237999999999
$ws.Range("I132").Value = 2893.3333
$ws.Range("K132").Value = 8679.999899999999
$ws.Range("M132").Value = -6149.999899999999
$ws.Range("H136").Value = 20842376
$ws.Range("I136").Value = 71430820
$ws.Range("J136").Value = 11843.471
$ws.Range("K136").Value = 214292460
$ws.Range("L136").Value = 35530.413
$ws.Range("M136").Value = -214289910
$ws.Range("N136").Value = -40630.413
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 47.444443
$ws.Range("I38").Value = 48.666668
$ws.Range("K38").Value = 146.000004
$ws.Range("M38").Value = 200.999996
$ws.Range("H107").Value = 13333629
$ws.Range("J107").Value = 18182090
$ws.Range("L107").Value = 54546270
$ws.Range("N107").Value = -54550110
$ws.Range("H121").Value = 1621
$ws.Range("J121").Value = 1985.1666
$ws.Range("L121").Value = 5955.4998
$ws.Range("N121").Value = -8575.4998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5918.5
$ws.Range("I113").Value = 2609
$ws.Range("J113").Value = 8124.8335
$ws.Range("K113").Value = 2609
$ws.Range("L113").Value = 8124.8335
$ws.Range("M113").Value = -439
$ws.Range("N113").Value = -12464.8335
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1320
$ws.Range("I16").Value = 810.25
$ws.Range("J16").Value = 1999.6666
$ws.Range("K16").Value = 810.25
$ws.Range("L16").Value = 1999.6666
$ws.Range("M16").Value = -640.25
$ws.Range("N16").Value = -2339.6666
$ws.Range("H40").Value = 5162.2607
$ws.Range("I40").Value = 3609.7693
$ws.Range("K40").Value = 3609.7693
$ws.Range("M40").Value = -3473.7693
$ws.Range("H46").Value = 1867
$ws.Range("J46").Value = 2904.25
$ws.Range("L46").Value = 2904.25
$ws.Range("N46").Value = -3280.25
$ws.Range("H68").Value = 5499.75
$ws.Range("J68").Value = 6333
$ws.Range("L68").Value = 6333
$ws.Range("N68").Value = -7831
$ws.Range("H71").Value = 5499.75
$ws.Range("J71").Value = 6333
$ws.Range("L71").Value = 31665
$ws.Range("N71").Value = -39153
$ws.Range("H82").Value = 641412.9
$ws.Range("I82").Value = 939825.7
$ws.Range("J82").Value = 1956.8572
$ws.Range("K82").Value = 939825.7
$ws.Range("L82").Value = 1956.8572
$ws.Range("M82").Value = -939464.7
$ws.Range("N82").Value = -2678.8572
$ws.Range("H85").Value = 641412.9
$ws.Range("I85").Value = 939825.7
$ws.Range("J85").Value = 1956.8572
$ws.Range("K85").Value = 939825.7
$ws.Range("L85").Value = 1956.8572
$ws.Range("M85").Value = -938577.7
$ws.Range("N85").Value = -4452.8572
$ws.Range("H122").Value = 4993.615
$ws.Range("I122").Value = 2958.6428
$ws.Range("J122").Value = 7367.75
$ws.Range("K122").Value = 8875.928400000001
$ws.Range("L122").Value = 22103.25
$ws.Range("M122").Value = -6425.928400000001
$ws.Range("N122").Value = -27003.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1833.3334
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 1833.3334
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H96").Value = 3639.111
$ws.Range("J96").Value = 1404
$ws.Range("L96").Value = 1404
$ws.Range("N96").Value = -4150
$ws.Range("H122").Value = 4010.7593
$ws.Range("J122").Value = 5115
$ws.Range("L122").Value = 15345
$ws.Range("N122").Value = -20245
$ws.Range("H126").Value = 3465.0588
$ws.Range("I126").Value = 2377
$ws.Range("K126").Value = 7131
$ws.Range("M126").Value = -4661
$ws.Range("H136").Value = 40006950
$ws.Range("I136").Value = 142858540
$ws.Range("K136").Value = 428575620
$ws.Range("M136").Value = -428573070
